$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.175.60"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "'1.872.17"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'307.52"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.5147"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "'0.3761"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "'0.07180"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "'0.8885"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "'20.72"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "'0.07583"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").Value = "'1.854.41"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "'5.343"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "'89.44"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'0.000008561"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "'14.19"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'27.223.99"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").Value = "'5.070"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Value = "'2.095.79"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "'6.498"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").Value = "'151.04"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").Value = "'1.846"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").Value = "'18.03"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "'2.135"
$ws.Range("E28").Value = "  -3.87%  "
$ws.Range("D29").Value = "'112.83"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").Value = "'4.765"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").Value = "'4.698"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'0.09003"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'0.05151"
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").Value = "'3.104"
$ws.Range("E34").Value = "  -3.27%  "
$ws.Range("D35").Value = "'0.7547"
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("D36").Value = "'1.174"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").Value = "'0.02044"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "'2.530"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "'3.026"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'1.080"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").Value = "'0.5366"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").Value = "'6.650"
$ws.Range("E42").Value = "  -4.26%  "
$ws.Range("D43").Value = "'114.74"
$ws.Range("E43").Value = "  +3.26%  "
$ws.Range("D44").Value = "'8.558"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "'0.1484"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "'0.4681"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'10.15"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").Value = "'1.575"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").Value = "'65.13"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").Value = "'36.53"
$ws.Range("E51").Value = "  -1.27%  "
